$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("area_pop_sum")

# Rename "Population" -> "population"
$ws.Range("A3").Value = "population"

# Move the density value from column C (C2/C3) into a new row 4
$ws.Range("A4").Value = "density"
$ws.Range("B4").Value = 1140.11776772289

# Clear out column C entirely (C1 "Density", C2 and C3 values)
$ws.Range("C1:C3").Clear()
